$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1199.746215820312
$ws.Range("C2").Value = 0.9
$ws.Range("D2").Value = 0.9009000062942505
$ws.Range("E2").Value = 1.228800058364868
$ws.Range("F2").Value = 0.5206000208854675
$ws.Range("H2").Value = 0.5258

# Row 3
$ws.Range("B3").Value = 1200.473999023438
$ws.Range("C3").Value = 0.9604
$ws.Range("D3").Value = 0.9451000000000001
$ws.Range("E3").Value = 1.291599988937378
$ws.Range("F3").Value = 0.6414999961853027
$ws.Range("H3").Value = 0.9181

# Row 4
$ws.Range("B4").Value = 811.87109375
$ws.Range("C4").Value = 0.9484
$ws.Range("D4").Value = 0.9252
$ws.Range("E4").Value = 1.420799970626831
$ws.Range("F4").Value = 0.7504000067710876
$ws.Range("H4").Value = 0.7417

# Row 5
$ws.Range("B5").Value = 783.7109985351562
$ws.Range("C5").Value = 0.8267
$ws.Range("D5").Value = 0.8232
$ws.Range("E5").Value = 1.062299966812134
$ws.Range("F5").Value = 0.4384999871253967
$ws.Range("H5").Value = -0.1619

# Row 6
$ws.Range("B6").Value = 1091.20751953125
$ws.Range("C6").Value = 0.8633
$ws.Range("D6").Value = 0.8668
$ws.Range("E6").Value = 1.053599953651428
$ws.Range("F6").Value = 0.5594000220298767
$ws.Range("H6").Value = 0.2237

# Row 7
$ws.Range("B7").Value = 861.7249145507812
$ws.Range("C7").Value = 0.8678
$ws.Range("D7").Value = 0.8694999814033508
$ws.Range("E7").Value = 1.012400031089783
$ws.Range("F7").Value = 0.6894999742507935
$ws.Range("H7").Value = 0.2478

# Row 8
$ws.Range("B8").Value = 954.5521240234375
$ws.Range("C8").Value = 0.8552999999999999
$ws.Range("D8").Value = 0.8522
$ws.Range("E8").Value = 1.045899987220764
$ws.Range("F8").Value = 0.708299994468689
$ws.Range("H8").Value = 0.0946

# Row 9
$ws.Range("B9").Value = 6903.287109375
$ws.Range("C9").Value = 0.8895999999999999
$ws.Range("D9").Value = 0.8848
$ws.Range("E9").Value = 1.420799970626831
$ws.Range("F9").Value = 0.4384999871253967
$ws.Range("H9").Value = 2.589799999999999
